$d = $word.ActiveDocument

$d.Content.Find.Execute("326÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "447÷4=", 2) | Out-Null
$d.Content.Find.Execute("481÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷6=", 2) | Out-Null
$d.Content.Find.Execute("741÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷8=", 2) | Out-Null
$d.Content.Find.Execute("272÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "493÷6=", 2) | Out-Null
$d.Content.Find.Execute("529÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "538÷5=", 2) | Out-Null
$d.Content.Find.Execute("972÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷9=", 2) | Out-Null
$d.Content.Find.Execute("342÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "373÷3=", 2) | Out-Null
$d.Content.Find.Execute("978÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "491÷7=", 2) | Out-Null
$d.Content.Find.Execute("688÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "312÷3=", 2) | Out-Null
$d.Content.Find.Execute("743÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "467÷2=", 2) | Out-Null
$d.Content.Find.Execute("341÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "172÷8=", 2) | Out-Null
$d.Content.Find.Execute("206÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "999÷2=", 2) | Out-Null
$d.Content.Find.Execute("453÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "937÷7=", 2) | Out-Null
$d.Content.Find.Execute("631÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "918÷9=", 2) | Out-Null
$d.Content.Find.Execute("875÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷3=", 2) | Out-Null
$d.Content.Find.Execute("481÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "288÷5=", 2) | Out-Null
$d.Content.Find.Execute("602÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷7=", 2) | Out-Null
$d.Content.Find.Execute("171÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "688÷3=", 2) | Out-Null
$d.Content.Find.Execute("831÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "298÷5=", 2) | Out-Null
$d.Content.Find.Execute("218÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "784÷5=", 2) | Out-Null
$d.Content.Find.Execute("916÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "864÷4=", 2) | Out-Null
$d.Content.Find.Execute("200÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷6=", 2) | Out-Null
$d.Content.Find.Execute("895÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "774÷5=", 2) | Out-Null
$d.Content.Find.Execute("451÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "891÷2=", 2) | Out-Null
$d.Content.Find.Execute("750÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "938÷5=", 2) | Out-Null
